$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G
    $val = $cell.Value2
    if ($val -eq $null) { continue }

    $parts = $val -split ', '
    $systemParts = @()
    $otherParts = @()
    foreach ($p in $parts) {
        if ($p -eq 'System' -or $p -eq 'system') {
            $systemParts += $p
        } else {
            $otherParts += $p
        }
    }

    if ($systemParts.Count -gt 0 -and $otherParts.Count -gt 0) {
        $newVal = ($otherParts + $systemParts) -join ', '
        if ($newVal -ne $val) {
            $cell.Value2 = $newVal
        }
    }
}
